$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Planning prévisionnel")
$ws1.Activate()
$excel.ActiveWindow.Zoom = 55
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("D30").Select()
